$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("B8").Value = "OKTP"
$ws.Range("D8").Value = "Ok to Plate"
$ws.Range("L8").Value = ""
$ws.Range("O8").Value = "OKTP"
$ws.Range("P8").Value = "Ok to Plate`nIntegration Validation"

# Row 9
$ws.Range("B9").Value = "WF Print F 4x0"
$ws.Range("D9").Value = "Placeholder - Integration Dal"
$ws.Range("L9").Value = "ISM Dallas"
$ws.Range("O9").Value = "Integration Validation"
$ws.Range("P9").Value = "Placeholder - Integration Dal"

# Row 10
$ws.Range("B10").Value = "Cut"
$ws.Range("D10").Value = "Zund 3XL3200 - #1"
$ws.Range("O10").Value = "Zund 3XL3200 - #1"
$ws.Range("P10").Value = "Zund 3XL3200 - #1"
